$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "summ00779629"
$wb.Worksheets.Item(2).Name = "summ00982942"
$wb.Worksheets.Item(3).Name = "summ01181446"
$wb.Worksheets.Item(4).Name = "summ01390430"
$wb.Worksheets.Item(5).Name = "summ01583214"
$wb.Worksheets.Item(6).Name = "summ01781234"
$wb.Worksheets.Item(7).Name = "summ01979258"
$wb.Worksheets.Item(8).Name = "summ02170790"
$wb.Worksheets.Item(9).Name = "summ02368428"
$wb.Worksheets.Item(10).Name = "summ02567621"
$wb.Worksheets.Item(11).Name = "summ02799271"
$wb.Worksheets.Item(12).Name = "summ03001292"
$wb.Worksheets.Item(13).Name = "summ03333641"
$wb.Worksheets.Item(14).Name = "summ03548234"
$wb.Worksheets.Item(15).Name = "summ03754901"
$wb.Worksheets.Item(16).Name = "summ03950044"
$wb.Worksheets.Item(17).Name = "summ04148065"
$wb.Worksheets.Item(18).Name = "summ04341092"
$wb.Worksheets.Item(19).Name = "summ04538626"
$wb.Worksheets.Item(20).Name = "summ04738776"
$wb.Worksheets.Item(21).Name = "summ04963711"
$wb.Worksheets.Item(22).Name = "summ05158507"
$wb.Worksheets.Item(23).Name = "summ05356020"
$wb.Worksheets.Item(24).Name = "summ05549040"
$wb.Worksheets.Item(25).Name = "summ05749367"
$wb.Worksheets.Item(26).Name = "summ05944516"
$wb.Worksheets.Item(27).Name = "summ06136536"
$wb.Worksheets.Item(28).Name = "summ06331567"
$wb.Worksheets.Item(29).Name = "summ06530590"
$wb.Worksheets.Item(30).Name = "summ06723676"
$wb.Worksheets.Item(31).Name = "summ06924700"
$wb.Worksheets.Item(32).Name = "summ07118468"
$wb.Worksheets.Item(33).Name = "summ07314485"
$wb.Worksheets.Item(34).Name = "summ07510180"
$wb.Worksheets.Item(35).Name = "summ07703253"
$wb.Worksheets.Item(36).Name = "summ07899276"
$wb.Worksheets.Item(37).Name = "summ08095294"
$wb.Worksheets.Item(38).Name = "summ08292318"
$wb.Worksheets.Item(39).Name = "summ08482337"
$wb.Worksheets.Item(40).Name = "summ08680362"
$wb.Worksheets.Item(41).Name = "summ08873381"
$wb.Worksheets.Item(42).Name = "summ09072924"
$wb.Worksheets.Item(43).Name = "summ09275440"
$wb.Worksheets.Item(44).Name = "summ09472462"
$wb.Worksheets.Item(45).Name = "summ09667484"
$wb.Worksheets.Item(46).Name = "summ09866743"
$wb.Worksheets.Item(47).Name = "summ10068766"
$wb.Worksheets.Item(48).Name = "summ10259785"
$wb.Worksheets.Item(49).Name = "summ10458164"
$wb.Worksheets.Item(50).Name = "summ10655721"
